$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "43.695.48"
Set-TextValue "E2" "  -0.30%  "
Set-TextValue "D3" "2.251.99"
Set-TextValue "E3" "  -1.01%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "232.08"
Set-TextValue "E5" "  -0.19%  "
Set-TextValue "E6" "  +1.44%  "
Set-TextValue "D7" "63.36"
Set-TextValue "E7" "  -1.68%  "
Set-TextValue "E8" "  -0.04%  "
Set-TextValue "D9" "0.445"
Set-TextValue "E9" "  +3.43%  "
Set-TextValue "D10" "0.0968"
Set-TextValue "E10" "  -5.98%  "
Set-TextValue "D11" "57.47"
Set-TextValue "E11" "  -0.07%  "
Set-TextValue "D12" "26.37"
Set-TextValue "E12" "  +0.50%  "
Set-TextValue "E13" "  +1.59%  "
Set-TextValue "D14" "2.586.40"
Set-TextValue "E14" "  -0.89%  "
Set-TextValue "D15" "15.53"
Set-TextValue "E15" "  -1.29%  "
Set-TextValue "E16" "  +2.08%  "
Set-TextValue "D17" "0.838"
Set-TextValue "E17" "  +1.36%  "
Set-TextValue "D18" "2.260.64"
Set-TextValue "E18" "  -0.20%  "
Set-TextValue "D19" "43.651.15"
Set-TextValue "E19" "  -0.11%  "
Set-TextValue "D20" "0.0₃0973"
Set-TextValue "E20" "  -2.72%  "
Set-TextValue "D21" "73.45"
Set-TextValue "E21" "  -0.75%  "
Set-TextValue "E22" "  +0.20%  "
Set-TextValue "D23" "247.75"
Set-TextValue "E23" "  -1.16%  "
Set-TextValue "D24" "1.00"
Set-TextValue "E24" "  -0.05%  "
Set-TextValue "D25" "3.63"
Set-TextValue "E25" "  +29.90%  "
Set-TextValue "E26" "  -2.53%  "
Set-TextValue "D27" "2.30"
Set-TextValue "E27" "  -1.12%  "
Set-TextValue "D28" "9.87"
Set-TextValue "E28" "  -1.10%  "
Set-TextValue "D29" "173.30"
Set-TextValue "E29" "  -0.02%  "
Set-TextValue "D30" "21.72"
Set-TextValue "E30" "  +3.30%  "
Set-TextValue "D31" "0.135"
Set-TextValue "E31" "  -1.11%  "
Set-TextValue "E32" "  -1.79%  "
Set-TextValue "E33" "  +1.46%  "
Set-TextValue "D34" "4.91"
Set-TextValue "E34" "  +3.55%  "
Set-TextValue "E35" "  -1.40%  "
Set-TextValue "E36" "  -3.02%  "
Set-TextValue "D37" "3.66"
Set-TextValue "E37" "  -4.44%  "
Set-TextValue "D38" "6.35"
Set-TextValue "E38" "  -5.86%  "
Set-TextValue "D39" "2.29"
Set-TextValue "E39" "  -2.65%  "
Set-TextValue "E40" "  +1.40%  "
Set-TextValue "E41" "  +0.19%  "
Set-TextValue "D42" "8.60"
Set-TextValue "E42" "  +2.33%  "
Set-TextValue "B43" "FTXToken"
Set-TextValue "C43" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D43" "4.65"
Set-TextValue "E43" "  +3.49%  "
Set-TextValue "D44" "98.03"
Set-TextValue "E44" "  +0.01%  "
Set-TextValue "B45" "InjectiveProtocol"
Set-TextValue "C45" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "17.11"
Set-TextValue "E45" "  -3.31%  "
Set-TextValue "D46" "0.0943"
Set-TextValue "E46" "  -2.79%  "
Set-TextValue "E47" "  -1.96%  "
Set-TextValue "D48" "1.452.85"
Set-TextValue "E48" "  -1.85%  "
Set-TextValue "E49" "  +1.55%  "
Set-TextValue "E50" "  -2.53%  "
Set-TextValue "D51" "9.79"
Set-TextValue "E51" "  -6.86%  "

"Applied 88 cell updates"
